# TC12_CDS_Filter_Study-WashingtonUniversity.xlsx : "CDS Study filter fixes"
#
# The ParticipantsTab query cell (B2 on the "startup" sheet) is replaced
# with an updated Neo4j/Cypher query (sorted samples, extra OPTIONAL MATCH
# clauses, reformatted RETURN/ORDER BY/LIMIT), and the active selection
# moves from A2 to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Washington University PDX Development and Trial Center"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

# @"..."@ here-strings keep a trailing newline before the closing tag;
# strip it so the cell content matches the source query exactly.
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newQuery

$ws.Range("C3").Select()
